$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.643.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.532.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3928"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3168"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.052"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.692"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.591"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.543.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06600"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.108"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.348"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.645.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.334"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.844"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.717.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.028"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9272"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -15.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.501"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.153"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02220"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.453"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2024"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.180"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5767"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5499"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.19%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.874"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06679"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.00%  "
